$wb = $excel.ActiveWorkbook

# Sheet "展览" (rId1 / sheet1.xml)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 9407
$ws1.Range("F6").Value = 9407
$ws1.Range("F7").Value = 596
$ws1.Range("F8").Value = 101
$ws1.Range("F9").Value = 157
$ws1.Range("F10").Value = 273
$ws1.Range("F11").Value = 406
$ws1.Range("F12").Value = 151
$ws1.Range("F13").Value = 179
$ws1.Range("F14").Value = 433
$ws1.Range("F15").Value = 12027
$ws1.Range("F16").Value = 12027
$ws1.Range("F20").Value = 32
$ws1.Range("F22").Value = 153
$ws1.Range("F24").Value = 237
$ws1.Range("F29").Value = 2723
$ws1.Range("F33").Value = 67
$ws1.Range("F34").Value = 14
$ws1.Range("F37").Value = 999
$ws1.Range("F38").Value = 4191
$ws1.Range("F39").Value = 3628
$ws1.Range("F40").Value = 526
$ws1.Range("F43").Value = 1317
$ws1.Range("F45").Value = 773
$ws1.Range("F46").Value = 416
$ws1.Range("F47").Value = 512
$ws1.Range("F49").Value = 216
$ws1.Range("F50").Value = 126
$ws1.Range("F51").Value = 135

# Sheet "演出" (rId2 / sheet2.xml)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F17").Value = 9
$ws2.Range("F19").Value = 9
$ws2.Range("F21").Value = 4

# Sheet "全部类型" (rId4 / sheet4.xml)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 9407
$ws4.Range("F10").Value = 596
$ws4.Range("F11").Value = 101
$ws4.Range("F12").Value = 157
$ws4.Range("F13").Value = 273
$ws4.Range("F14").Value = 406
$ws4.Range("F15").Value = 151
$ws4.Range("F16").Value = 179
$ws4.Range("F17").Value = 12027
$ws4.Range("F18").Value = 12027
$ws4.Range("F20").Value = 32
$ws4.Range("F24").Value = 153
$ws4.Range("F30").Value = 2723
$ws4.Range("F34").Value = 67
$ws4.Range("F35").Value = 14
$ws4.Range("F37").Value = 9
$ws4.Range("F39").Value = 999
$ws4.Range("F42").Value = 3628
$ws4.Range("F45").Value = 1317
$ws4.Range("F47").Value = 416
$ws4.Range("F49").Value = 512
$ws4.Range("F51").Value = 216
